# Applies the 2026-01-12 export changes:
#  - Fix a few employee/client names on the "Weekly Timesheet" sheet
#  - Fill in the simulator's Rate/Total numbers (previously all zero) on both sheets
#  - Update the Employee ID shared string

$wb = $excel.ActiveWorkbook

$wsTime   = $wb.Worksheets.Item("Weekly Timesheet")
$wsSchema = $wb.Worksheets.Item("Jason Schema")

# --- Fix client names (shared string correction touches both sheets,
#     since "Weekly Timesheet"!B and "Jason Schema"!D mirror the same
#     underlying client names) ---
$wsTime.Range("B3").Value = "Evans"    # was "Richer"
$wsTime.Range("B5").Value = "Hewett"   # was "Tercek"
$wsTime.Range("B6").Value = "Howard"   # was "Patton"

$wsSchema.Range("D3").Value = "Evans"    # was "Richer"
$wsSchema.Range("D5").Value = "Hewett"   # was "Tercek"
$wsSchema.Range("D6").Value = "Howard"   # was "Patton"

# --- Weekly Timesheet: Rate (E) / Total (F) for each daily row ---
$rate  = 140
$total = 1120

foreach ($r in 2..6) {
    $wsTime.Cells.Item($r, 5).Value = $rate   # column E - Rate
    $wsTime.Cells.Item($r, 6).Value = $total  # column F - Total
}

# --- Weekly Timesheet: subtotal / grand total rows ---
$wsTime.Range("F8").Value  = 5600   # HOURLY SUBTOTAL
$wsTime.Range("F12").Value = 5600   # ADMIN SUBTOTAL
$wsTime.Range("F13").Value = 5600   # GRAND TOTAL

# --- Jason Schema: Rate (F) / Total (G) for each daily row ---
foreach ($r in 2..6) {
    $wsSchema.Cells.Item($r, 6).Value = $rate   # column F - Rate
    $wsSchema.Cells.Item($r, 7).Value = $total  # column G - Total
}

# --- Employee ID correction (shared string used on Jason Schema sheet) ---
$employeeId = "emp_jp4mlvog"   # was "emp_5chpvt65"
foreach ($r in 2..6) {
    $wsSchema.Cells.Item($r, 2).Value = $employeeId
}
